# Slide 27 has the console-output textbox (Shape 635) reading:
#   " $ python tryexcept.py "
#   " First -1"
#   " Second 123"
# Sue bolds that whole textbox (one run in it was already bold; this
# brings the remaining runs to b="1" as well).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(27)
$shp = $s.Shapes.Item("Shape 635")
$shp.TextFrame.TextRange.Font.Bold = -1
